# fix: ubah tahun ke 2026
#
# The template referenced "tahun 2024" (twice) and "Januari 2024" (once).
# Bump every occurrence of the year from 2024 to 2026.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$replaced = $find.Execute(
    "2024",   # FindText
    $false,   # MatchCase
    $false,   # MatchWholeWord
    $false,   # MatchWildcards
    $false,   # MatchSoundsLike
    $false,   # MatchAllWordForms
    $true,    # Forward
    1,        # Wrap (wdFindContinue)
    $false,   # Format
    "2026",   # ReplaceWith
    2         # Replace (wdReplaceAll)
)

if (-not $replaced) {
    throw "Could not find '2024' in the document to replace with '2026'."
}
